# Insert two new translation rows into Sheet1 ("CriterionFormula" / "CalculatedCrit")
# right after the existing "cf_min_Magnesium_mg/l" row (row 80), pushing the
# "IR_ActivityType" row (and everything after it) down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank rows before row 81 (old row 81 = "IR_ActivityType")
$ws.Rows.Item(81).Resize(2).Insert()

$ws.Range("A81").Value = "CriterionFormula"
$ws.Range("B81").Value = "DA"
$ws.Range("A82").Value = "CalculatedCrit"
$ws.Range("B82").Value = "DA"
